$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Element Name:) ---
$ws.Range("F2").Value2 = "flight-departing-txt-id"
$ws.Range("G2").Value2 = "flight-returning-txt-id"
$ws.Range("H2").Value2 = "flight-add-hotel-ckb-id"
$ws.Range("I2").Value2 = "flight-hotel-checkin-txt-id"
$ws.Range("J2").Value2 = "flight-hotel-checkout-txt-id"
$ws.Range("K2").Value2 = "search-btn-class"
$ws.Range("L2").Value2 = "search-auto-id"

# --- Row 3 (Element Value:) ---
$ws.Range("F3").Value2 = "flight-departing-hp-flight"
$ws.Range("G3").Value2 = "flight-returning-hp-flight"
$ws.Range("H3").Value2 = "flight-add-hotel-checkbox-hp-flight"
$ws.Range("I3").Value2 = "flight-hotel-checkin-hp-flight"
$ws.Range("J3").Value2 = "flight-hotel-checkout-hp-flight"
$ws.Range("K3").Value2 = "gcw-submit"
$ws.Range("L3").Value2 = "aria-option-0"

# --- Column widths to match new layout ---
$ws.Columns.Item(6).ColumnWidth = 18.833333333333332
$ws.Columns.Item(7).ColumnWidth = 18.833333333333332
$ws.Columns.Item(8).ColumnWidth = 25.333333333333332
$ws.Columns.Item(9).ColumnWidth = 20.666666666666668
$ws.Columns.Item(10).ColumnWidth = 21.666666666666668
$ws.Columns.Item(11).ColumnWidth = 12.5
$ws.Columns.Item(12).ColumnWidth = 10.833333333333334

# --- Selection moves to I18 ---
$ws.Range("I18").Select()
